$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.330.09'
$ws.Range('E2').Value = '  -2.43%  '

$ws.Range('D3').Value = '1.933.86'
$ws.Range('E3').Value = '  -2.32%  '

$ws.Range('E4').Value = '  -0.47%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.69'
$ws.Range('E5').Value = '  -1.40%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.7090'
$ws.Range('E6').Value = '  -3.58%  '

$ws.Range('E7').Value = '  -0.46%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3302'
$ws.Range('E8').Value = '  -2.53%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '27.67'
$ws.Range('E9').Value = '  +0.28%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07302'
$ws.Range('E10').Value = '  +2.10%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8052'
$ws.Range('E11').Value = '  -2.89%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08083'
$ws.Range('E12').Value = '  -0.22%  '

$ws.Range('D13').Value = '1.933.49'
$ws.Range('E13').Value = '  -2.36%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.474'
$ws.Range('E14').Value = '  -2.25%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.60'
$ws.Range('E15').Value = '  -4.71%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.10'
$ws.Range('E16').Value = '  -1.86%  '

$ws.Range('D17').Value = '30.329.61'
$ws.Range('E17').Value = '  -2.46%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '253.40'
$ws.Range('E18').Value = '  -5.40%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008184'
$ws.Range('E19').Value = '  -0.83%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.791'
$ws.Range('E20').Value = '  -4.71%  '

$ws.Range('D21').Value = '2.189.68'
$ws.Range('E21').Value = '  -2.80%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  -0.49%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('E23').Value = '  -0.37%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.977'
$ws.Range('E24').Value = '  -1.23%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.749'
$ws.Range('E25').Value = '  -2.37%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.34'
$ws.Range('E26').Value = '  +2.08%  '

$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.348'
$ws.Range('E27').Value = '  -0.40%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.33'
$ws.Range('E28').Value = '  -1.98%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1284'
$ws.Range('E29').Value = '  -3.29%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.349'
$ws.Range('E30').Value = '  -2.27%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.541'
$ws.Range('E31').Value = '  -3.34%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.418'
$ws.Range('E32').Value = '  -4.58%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.174'
$ws.Range('E33').Value = '  -5.42%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05199'
$ws.Range('E34').Value = '  -1.92%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.263'
$ws.Range('E35').Value = '  -1.28%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7468'
$ws.Range('E36').Value = '  -4.29%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.787'
$ws.Range('E37').Value = '  -0.12%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01966'
$ws.Range('E38').Value = '  -1.92%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.808'
$ws.Range('E39').Value = '  -2.43%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '79.00'

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.416'
$ws.Range('E41').Value = '  -5.29%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4524'
$ws.Range('E42').Value = '  -2.46%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.016'
$ws.Range('E43').Value = '  -4.18%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8455'
$ws.Range('E44').Value = '  -0.98%  '

$ws.Range('E45').Value = '  -0.49%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.57'
$ws.Range('E46').Value = '  -3.04%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.753'
$ws.Range('E47').Value = '  -3.11%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.451'
$ws.Range('E48').Value = '  -2.57%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '36.74'
$ws.Range('E49').Value = '  -1.43%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4171'
$ws.Range('E50').Value = '  -3.30%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06034'
$ws.Range('E51').Value = '  +0.02%  '
